# Applies the scheduled market-data refresh to the Atomos_Profits workbook.
# For each touched leve row, currentAveragePrice / NQ / HQ columns (H-N)
# are updated to the latest pulled values; profit columns (M/N) are
# recomputed accordingly. Only cells whose values actually moved are
# written, so pre-existing gaps (e.g. a missing LeveProfitNQ cell) are
# preserved unless the diff explicitly adds/removes them.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1648.2222
$ws.Range("I43").Value = 995.6667
$ws.Range("J43").Value = 1778.7333
$ws.Range("K43").Value = 995.6667
$ws.Range("L43").Value = 1778.7333
$ws.Range("M43").Value = -926.6667
$ws.Range("N43").Value = -1916.7333

$ws.Range("H125").Value = 1387.6666
$ws.Range("I125").Value = 1430.1818
$ws.Range("J125").Value = 1351.6923
$ws.Range("K125").Value = 12871.6362
$ws.Range("L125").Value = 12165.2307
$ws.Range("M125").Value = -10411.6362
$ws.Range("N125").Value = -17085.2307

$ws.Range("H135").Value = 1231.5
$ws.Range("I135").Value = 1035.4667
$ws.Range("J135").Value = 1819.6
$ws.Range("K135").Value = 9319.2003
$ws.Range("L135").Value = 16376.4
$ws.Range("M135").Value = -6784.2003
$ws.Range("N135").Value = -21446.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2059.35
$ws.Range("I2").Value = 1743
$ws.Range("J2").Value = 2318.182
$ws.Range("K2").Value = 1743
$ws.Range("L2").Value = 2318.182
$ws.Range("M2").Value = -1630
$ws.Range("N2").Value = -2544.182

$ws.Range("H116").Value = 2059.35
$ws.Range("I116").Value = 1743
$ws.Range("J116").Value = 2318.182
$ws.Range("K116").Value = 1743
$ws.Range("L116").Value = 2318.182
$ws.Range("M116").Value = 551
$ws.Range("N116").Value = -6906.182

$ws.Range("H132").Value = 1404.75
$ws.Range("I132").Value = 1161.1666
$ws.Range("J132").Value = 3012.4
$ws.Range("K132").Value = 3483.4998
$ws.Range("L132").Value = 9037.200000000001
$ws.Range("M132").Value = -953.4998000000001
$ws.Range("N132").Value = -14097.2

$ws.Range("H135").Value = 26301.777
$ws.Range("J135").Value = 26301.777
$ws.Range("L135").Value = 26301.777
$ws.Range("N135").Value = -36441.777

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2059.35
$ws.Range("I3").Value = 1743
$ws.Range("J3").Value = 2318.182
$ws.Range("K3").Value = 1743
$ws.Range("L3").Value = 2318.182
$ws.Range("M3").Value = -1629
$ws.Range("N3").Value = -2546.182

$ws.Range("H99").Value = 2200.76
$ws.Range("I99").Value = 1554.0588
$ws.Range("J99").Value = 3575
$ws.Range("K99").Value = 1554.0588
$ws.Range("L99").Value = 3575
$ws.Range("M99").Value = -56.05880000000002
$ws.Range("N99").Value = -6571

$ws.Range("H134").Value = 1607.807
$ws.Range("I134").Value = 1288.5834
$ws.Range("J134").Value = 3310.3333
$ws.Range("K134").Value = 3865.7502
$ws.Range("L134").Value = 9930.999899999999
$ws.Range("M134").Value = -1330.7502
$ws.Range("N134").Value = -15000.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 10419598
$ws.Range("I58").Value = 2118
$ws.Range("J58").Value = 31254558
$ws.Range("K58").Value = 2118
$ws.Range("L58").Value = 31254558
$ws.Range("M58").Value = -1915
$ws.Range("N58").Value = -31254964

$ws.Range("H94").Value = 4827.077
$ws.Range("I94").Value = 11504
$ws.Range("J94").Value = 2824
$ws.Range("K94").Value = 11504
$ws.Range("L94").Value = 2824
$ws.Range("M94").Value = -11053
$ws.Range("N94").Value = -3726

$ws.Range("H107").Value = 1083.8148
$ws.Range("I107").Value = 384.6
$ws.Range("K107").Value = 384.6
$ws.Range("M107").Value = 1535.4

$ws.Range("H132").Value = 3105.8462
$ws.Range("I132").Value = 2520.7058
$ws.Range("J132").Value = 4211.1113
$ws.Range("K132").Value = 7562.117400000001
$ws.Range("L132").Value = 12633.3339
$ws.Range("M132").Value = -5032.117400000001
$ws.Range("N132").Value = -17693.3339

$ws.Range("H136").Value = 10419598
$ws.Range("I136").Value = 2118
$ws.Range("J136").Value = 31254558
$ws.Range("K136").Value = 6354
$ws.Range("L136").Value = 93763674
$ws.Range("M136").Value = -3804
$ws.Range("N136").Value = -93768774

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 3665.8333
$ws.Range("I76").Value = 1995
$ws.Range("K76").Value = 5985
$ws.Range("M76").Value = -5602

$ws.Range("H79").Value = 3665.8333
$ws.Range("I79").Value = 1995
$ws.Range("K79").Value = 5985
$ws.Range("M79").Value = -4659

$ws.Range("H80").Value = 3921.7778
$ws.Range("J80").Value = 4212
$ws.Range("L80").Value = 12636
$ws.Range("N80").Value = -14508

$ws.Range("H83").Value = 3921.7778
$ws.Range("J83").Value = 4212
$ws.Range("L83").Value = 37908
$ws.Range("N83").Value = -47268

$ws.Range("H133").Value = 5649.857
$ws.Range("I133").Value = 5926.5
$ws.Range("J133").Value = 3990
$ws.Range("K133").Value = 17779.5
$ws.Range("L133").Value = 11970
$ws.Range("M133").Value = -12719.5
$ws.Range("N133").Value = -22090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1517625
$ws.Range("I126").Value = 4546217
$ws.Range("J126").Value = 3329.05
$ws.Range("K126").Value = 13638651
$ws.Range("L126").Value = 9987.150000000001
$ws.Range("M126").Value = -13636181
$ws.Range("N126").Value = -14927.15

$ws.Range("H132").Value = 3359.795
$ws.Range("I132").Value = 3615.8096
$ws.Range("J132").Value = 3061.111
$ws.Range("K132").Value = 10847.4288
$ws.Range("L132").Value = 9183.332999999999
$ws.Range("M132").Value = -8317.4288
$ws.Range("N132").Value = -14243.333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125001600
$ws.Range("I22").Value = 142857540
$ws.Range("J22").Value = 10000
$ws.Range("K22").Value = 142857540
$ws.Range("L22").Value = 10000
$ws.Range("M22").Value = -142857245
$ws.Range("N22").Value = -10590

$ws.Range("H27").Value = 125001600
$ws.Range("I27").Value = 142857540
$ws.Range("J27").Value = 10000
$ws.Range("K27").Value = 142857540
$ws.Range("L27").Value = 10000
$ws.Range("M27").Value = -142857433
$ws.Range("N27").Value = -10214

$ws.Range("H46").Value = 1250.909
$ws.Range("I46").Value = 929.78723
$ws.Range("J46").Value = 3137.5
$ws.Range("K46").Value = 929.78723
$ws.Range("L46").Value = 3137.5
$ws.Range("M46").Value = -741.78723
$ws.Range("N46").Value = -3513.5

$ws.Range("H55").Value = 1087.2858
$ws.Range("I55").Value = 213.33333
$ws.Range("J55").Value = 1325.6364
$ws.Range("K55").Value = 213.33333
$ws.Range("L55").Value = 1325.6364
$ws.Range("M55").Value = -40.33332999999999
$ws.Range("N55").Value = -1671.6364

$ws.Range("H93").Value = 1823
$ws.Range("I93").Value = 701
$ws.Range("J93").Value = 2384
$ws.Range("K93").Value = 701
$ws.Range("L93").Value = 2384
$ws.Range("M93").Value = 547
$ws.Range("N93").Value = -4880

$ws.Range("H132").Value = 2358.1282
$ws.Range("I132").Value = 1699.7826
$ws.Range("J132").Value = 3304.5
$ws.Range("K132").Value = 5099.3478
$ws.Range("L132").Value = 9913.5
$ws.Range("M132").Value = -2569.3478
$ws.Range("N132").Value = -14973.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H10").Value = 80006
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 80006
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 80006
$ws.Range("N10").Value = -80344

$ws.Range("H15").Value = 24745.4
$ws.Range("I15").Value = 4006
$ws.Range("J15").Value = 29930.25
$ws.Range("K15").Value = 4006
$ws.Range("L15").Value = 29930.25
$ws.Range("M15").Value = -3718
$ws.Range("N15").Value = -30506.25

$ws.Range("H126").Value = 2943244.8
$ws.Range("I126").Value = 1443.9048
$ws.Range("J126").Value = 7695384.5
$ws.Range("K126").Value = 4331.7144
$ws.Range("L126").Value = 23086153.5
$ws.Range("M126").Value = -1861.7144
$ws.Range("N126").Value = -23091093.5

$ws.Range("H132").Value = 268597.66
$ws.Range("I132").Value = 359746.66
$ws.Range("J132").Value = 13380.5
$ws.Range("K132").Value = 1079239.98
$ws.Range("L132").Value = 40141.5
$ws.Range("M132").Value = -1076709.98
$ws.Range("N132").Value = -45201.5

$ws.Range("H136").Value = 1178.4054
$ws.Range("I136").Value = 555.2069
$ws.Range("J136").Value = 3437.5
$ws.Range("K136").Value = 1665.6207
$ws.Range("L136").Value = 10312.5
$ws.Range("M136").Value = 884.3793000000001
$ws.Range("N136").Value = -15412.5

# Row 10 on WVR no longer reports a LeveProfitNQ figure this refresh —
# the source cell is cleared so the column drops out for that row.
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M10").ClearContents()
